$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315056920051575
$ws.Range("B1").Value = 1.920963644981384
$ws.Range("C1").Value = 2.842329740524292
$ws.Range("D1").Value = 5.348701953887939
$ws.Range("E1").Value = 2.976485013961792
